$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

$names = @("GlobalPagesUser1","GlobalPagesUser2","GlobalPagesUser3","GlobalPagesUser4","GlobalPagesUser5","GlobalPagesUser6")
$emails = @("GlobalPagesUser1@mailinator.com","GlobalPagesUser2@mailinator.com","GlobalPagesUser3@mailinator.com","GlobalPagesUser4@mailinator.com","GlobalPagesUser5@mailinator.com","GlobalPagesUser6@mailinator.com")

$startRow = 114

# Column A (UserName) - fill all names first (matches shared-string insertion order)
for ($i = 0; $i -lt $names.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $names[$i]
}

# Column B (Password) - same password used by the rows above
for ($i = 0; $i -lt $names.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 2).Value = "Password1"
}

# Column E (Notes) + Column F (N) - copy formatting/values from the row above (row 113)
for ($i = 0; $i -lt $names.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 5).Value = "THIS IS IN USE 24/7 - DO NOT USE!"
    $ws.Cells.Item($r, 6).Value = "N"
    $ws.Range("E113:F113").Copy()
    $ws.Cells.Item($r, 5).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
}

# Column G (Email) - fill all emails (shared strings appended after the names, matching diff order)
for ($i = 0; $i -lt $emails.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 7).Value = $emails[$i]
}

# Add the mailto hyperlinks for column G, then restore the Hyperlink cell style
# (G113) so we don't pick up the engine's auto-generated duplicate style.
for ($i = 0; $i -lt $emails.Count; $i++) {
    $r = $startRow + $i
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 7), "mailto:" + $emails[$i])
    $ws.Range("G113").Copy()
    $ws.Cells.Item($r, 7).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
}

# Update the sheet view to match where Excel scrolled to/selected after the edit.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 94
$ws.Range("A120").Select()
